# v1: select palette color
#
# The original "select a palette color and apply to a region of the
# image / - fill mode / - undo, redo" task (Id 10, Active sheet row 5)
# is broken up into three separate tasks:
#   - "select a palette color"                         -> Id 31, marked Done, moved to Inactive sheet
#   - "fill in a section of color on the image"         -> Id 32, new Todo task on Active sheet
#   - "undo, redo coloring a section on the image"      -> Id 33, new Todo task on Active sheet
#
# The Max Id on the Config sheet is bumped from 30 to 33 to match.

$wb = $excel.ActiveWorkbook
$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")
$wsConfig = $wb.Worksheets.Item("Config")

# ---------------------------------------------------------------------
# 1) Inactive sheet: insert a new row 2 for the completed
#    "select a palette color" task (Id 31).
# ---------------------------------------------------------------------
$wsInactive.Rows.Item(2).Insert()
$wsInactive.Range("A2:F2").Style = "Normal"

$wsInactive.Cells.Item(2, 1).Value = 31
$wsInactive.Cells.Item(2, 2).Value = "select a palette color"
$wsInactive.Cells.Item(2, 3).Value = "Done"
$wsInactive.Cells.Item(2, 4).Value = "Task"

$wsInactive.Cells.Item(2, 5).NumberFormat = "@"
$wsInactive.Cells.Item(2, 5).Value = "8/11/2018"
$wsInactive.Cells.Item(2, 5).Style = "Normal"

$wsInactive.Cells.Item(2, 6).NumberFormat = "@"
$wsInactive.Cells.Item(2, 6).Value = "8/11/2018"
$wsInactive.Cells.Item(2, 6).Style = "Normal"

# ---------------------------------------------------------------------
# 2) Active sheet: replace the old row 5 ("select a palette color and
#    apply...") with two new todo rows: "fill in a section of color on
#    the image" (Id 32) and "undo, redo coloring a section on the
#    image" (Id 33).
# ---------------------------------------------------------------------
$wsActive.Rows.Item(5).Insert()
$wsActive.Range("A5:E5").Style = "Normal"

$wsActive.Cells.Item(5, 1).Value = 32
$wsActive.Cells.Item(5, 2).Value = "fill in a section of color on the image"
$wsActive.Cells.Item(5, 3).Value = "Todo"
$wsActive.Cells.Item(5, 4).Value = "Task"
$wsActive.Cells.Item(5, 5).NumberFormat = "@"
$wsActive.Cells.Item(5, 5).Value = "8/11/2018"
$wsActive.Cells.Item(5, 5).Style = "Normal"

# The former row 5 is now row 6; overwrite it with the second new task.
$wsActive.Cells.Item(6, 1).Value = 33
$wsActive.Cells.Item(6, 2).Value = "undo, redo coloring a section on the image"
$wsActive.Cells.Item(6, 3).Value = "Todo"
$wsActive.Cells.Item(6, 4).Value = "Task"
$wsActive.Cells.Item(6, 5).NumberFormat = "@"
$wsActive.Cells.Item(6, 5).Value = "8/11/2018"
$wsActive.Cells.Item(6, 5).Style = "Normal"

# ---------------------------------------------------------------------
# 3) Config sheet: bump the "Max Id" value from 30 to 33.
# ---------------------------------------------------------------------
$wsConfig.Cells.Item(2, 6).Value = 33
